# Docx writer: Use different style for block quotes in notes.
#
# Adds a new paragraph style "Footnote Block Text" (styleId
# "FootnoteBlockText"), based on "Footnote Text" (next paragraph style
# also "Footnote Text"), mirroring the existing "Block Text" style
# (which is based on "Body Text") so that block quotes that occur
# inside footnotes/endnotes can be styled independently from regular
# block quotes.

$d = $word.ActiveDocument

$style = $d.Styles.Add("Footnote Block Text", 1)   # 1 = wdStyleTypeParagraph

$style.BaseStyle           = "Footnote Text"
$style.NextParagraphStyle  = "Footnote Text"
$style.Priority            = 9
$style.UnhideWhenUsed      = $true
$style.QuickStyle          = $true

# Paragraph formatting: spacing before/after = 100 twips (5pt), and a
# left/right indent of 480 twips (24pt) with no first-line indent -
# same block-quote geometry "Block Text" uses.
$pf = $style.ParagraphFormat
$pf.SpaceBefore      = 5
$pf.SpaceAfter       = 5
$pf.FirstLineIndent  = 0
$pf.LeftIndent       = 24
$pf.RightIndent      = 24
